# Scheduled-runner style refresh of cached market-price / profit figures
# across the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only currentAveragePrice* / LevePrice* / LeveProfit* (columns H-N) cells
# are refreshed; identifying columns (A-G) are left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 6451
$ws.Range("I34").Value = 2715.7144
$ws.Range("K34").Value = 2715.7144
$ws.Range("M34").Value = -2512.7144
# Row 36
$ws.Range("H36").Value = 6451
$ws.Range("I36").Value = 2715.7144
$ws.Range("K36").Value = 2715.7144
$ws.Range("M36").Value = -2000.7144
# Row 40
$ws.Range("H40").Value = 1507.0555
$ws.Range("I40").Value = 1200
$ws.Range("J40").Value = 1660.5834
$ws.Range("K40").Value = 1200
$ws.Range("L40").Value = 1660.5834
$ws.Range("M40").Value = -1025
$ws.Range("N40").Value = -2010.5834
# Row 43
$ws.Range("H43").Value = 2254.2856
$ws.Range("I43").Value = 3812.5
$ws.Range("J43").Value = 1887.6471
$ws.Range("K43").Value = 3812.5
$ws.Range("L43").Value = 1887.6471
$ws.Range("M43").Value = -3743.5
$ws.Range("N43").Value = -2025.6471
# Row 63
$ws.Range("H63").Value = 19933.334
$ws.Range("J63").Value = 19933.334
$ws.Range("L63").Value = 19933.334
$ws.Range("N63").Value = -21181.334
# Row 66
$ws.Range("H66").Value = 19933.334
$ws.Range("J66").Value = 19933.334
$ws.Range("L66").Value = 59800.00199999999
$ws.Range("N66").Value = -66040.00199999999
# Row 112
$ws.Range("H112").Value = 1652.2609
$ws.Range("I112").Value = 750.3333
$ws.Range("J112").Value = 1970.5883
$ws.Range("K112").Value = 2250.9999
$ws.Range("L112").Value = 5911.7649
$ws.Range("M112").Value = -1142.9999
$ws.Range("N112").Value = -8127.7649
# Row 132
$ws.Range("H132").Value = 3675.568
$ws.Range("I132").Value = 2162.3076
$ws.Range("J132").Value = 5861.3887
$ws.Range("K132").Value = 6486.9228
$ws.Range("L132").Value = 17584.1661
$ws.Range("M132").Value = -3956.9228
$ws.Range("N132").Value = -22644.1661
# Row 138
$ws.Range("H138").Value = 1352.5411
$ws.Range("I138").Value = 795.7447
$ws.Range("K138").Value = 2387.2341
$ws.Range("M138").Value = 2752.7659
# Row 141
$ws.Range("H141").Value = 456.60938
$ws.Range("I141").Value = 452.0645
$ws.Range("J141").Value = 597.5
$ws.Range("K141").Value = 1356.1935
$ws.Range("L141").Value = 1792.5
$ws.Range("M141").Value = 3823.8065
$ws.Range("N141").Value = -12152.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 30177.5
$ws.Range("J24").Value = 30177.5
$ws.Range("L24").Value = 30177.5
$ws.Range("N24").Value = -30925.5
# Row 32
$ws.Range("H32").Value = 6108.467
$ws.Range("I32").Value = 4012.8975
$ws.Range("J32").Value = 19729.666
$ws.Range("K32").Value = 4012.8975
$ws.Range("L32").Value = 19729.666
$ws.Range("M32").Value = -3725.8975
$ws.Range("N32").Value = -20303.666
# Row 61
$ws.Range("H61").Value = 235133.83
$ws.Range("I61").Value = 158201.45
$ws.Range("J61").Value = 504397.1
$ws.Range("K61").Value = 158201.45
$ws.Range("L61").Value = 504397.1
$ws.Range("M61").Value = -157989.45
$ws.Range("N61").Value = -504821.1
# Row 74
$ws.Range("H74").Value = 854.99
$ws.Range("I74").Value = 643.25275
$ws.Range("J74").Value = 2995.889
$ws.Range("K74").Value = 643.25275
$ws.Range("L74").Value = 2995.889
$ws.Range("M74").Value = 230.74725
$ws.Range("N74").Value = -4743.889
# Row 77
$ws.Range("H77").Value = 854.99
$ws.Range("I77").Value = 643.25275
$ws.Range("J77").Value = 2995.889
$ws.Range("K77").Value = 3216.26375
$ws.Range("L77").Value = 14979.445
$ws.Range("M77").Value = 1151.73625
$ws.Range("N77").Value = -23715.445
# Row 96
$ws.Range("H96").Value = 18114.666
$ws.Range("J96").Value = 18114.666
$ws.Range("L96").Value = 18114.666
$ws.Range("N96").Value = -23606.666
# Row 100
$ws.Range("H100").Value = 30177.5
$ws.Range("J100").Value = 30177.5
$ws.Range("L100").Value = 30177.5
$ws.Range("N100").Value = -32341.5
# Row 136
$ws.Range("H136").Value = 235133.83
$ws.Range("I136").Value = 158201.45
$ws.Range("J136").Value = 504397.1
$ws.Range("K136").Value = 474604.35
$ws.Range("L136").Value = 1513191.3
$ws.Range("M136").Value = -472054.35
$ws.Range("N136").Value = -1518291.3

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 233.33333
$ws.Range("I80").Value = 97.333336
$ws.Range("J80").Value = 301.33334
$ws.Range("K80").Value = 97.333336
$ws.Range("L80").Value = 301.33334
$ws.Range("M80").Value = 900.666664
$ws.Range("N80").Value = -2297.33334
# Row 83
$ws.Range("H83").Value = 233.33333
$ws.Range("I83").Value = 97.333336
$ws.Range("J83").Value = 301.33334
$ws.Range("K83").Value = 486.66668
$ws.Range("L83").Value = 1506.6667
$ws.Range("M83").Value = 4505.33332
$ws.Range("N83").Value = -11490.6667
# Row 94
$ws.Range("H94").Value = 3749.75
$ws.Range("I94").Value = 1499.5
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 1499.5
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = -1048.5
$ws.Range("N94").Value = -6902

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2574.8225
$ws.Range("I31").Value = 1728.8572
$ws.Range("J31").Value = 4351.35
$ws.Range("K31").Value = 1728.8572
$ws.Range("L31").Value = 4351.35
$ws.Range("M31").Value = -1433.8572
$ws.Range("N31").Value = -4941.35
# Row 34
$ws.Range("H34").Value = 2574.8225
$ws.Range("I34").Value = 1728.8572
$ws.Range("J34").Value = 4351.35
$ws.Range("K34").Value = 1728.8572
$ws.Range("L34").Value = 4351.35
$ws.Range("M34").Value = -1526.8572
$ws.Range("N34").Value = -4755.35
# Row 58
$ws.Range("H58").Value = 2573.7354
$ws.Range("I58").Value = 2573.8147
$ws.Range("J58").Value = 2573.4285
$ws.Range("K58").Value = 2573.8147
$ws.Range("L58").Value = 2573.4285
$ws.Range("M58").Value = -2370.8147
$ws.Range("N58").Value = -2979.4285
# Row 87
$ws.Range("H87").Value = 15000
$ws.Range("J87").Value = 15000
$ws.Range("L87").Value = 15000
$ws.Range("N87").Value = -17372
# Row 90
$ws.Range("H90").Value = 15000
$ws.Range("J90").Value = 15000
$ws.Range("L90").Value = 45000
$ws.Range("N90").Value = -56856
# Row 112
$ws.Range("H112").Value = 26000
$ws.Range("J112").Value = 26000
$ws.Range("L112").Value = 26000
$ws.Range("N112").Value = -28954
# Row 118
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = $null
$ws.Range("N118").Value = 0
# Row 132
$ws.Range("H132").Value = 2230.5186
$ws.Range("I132").Value = 1307.2941
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 3921.8823
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -1391.8823
$ws.Range("N132").Value = -16460
# Row 134
$ws.Range("H134").Value = 1314.4762
$ws.Range("I134").Value = 787.875
$ws.Range("J134").Value = 2999.6
$ws.Range("K134").Value = 2363.625
$ws.Range("L134").Value = 8998.799999999999
$ws.Range("M134").Value = 171.375
$ws.Range("N134").Value = -14068.8
# Row 136
$ws.Range("H136").Value = 2573.7354
$ws.Range("I136").Value = 2573.8147
$ws.Range("J136").Value = 2573.4285
$ws.Range("K136").Value = 7721.4441
$ws.Range("L136").Value = 7720.2855
$ws.Range("M136").Value = -5171.4441
$ws.Range("N136").Value = -12820.2855

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 102
$ws.Range("H102").Value = 8666.666999999999
$ws.Range("J102").Value = 8875
$ws.Range("L102").Value = 26625
$ws.Range("N102").Value = -31493
# Row 113
$ws.Range("H113").Value = 576.3333
$ws.Range("J113").Value = 568.8333
$ws.Range("L113").Value = 1706.4999
$ws.Range("N113").Value = -6046.4999
# Row 130
$ws.Range("H130").Value = 202206.6
$ws.Range("J130").Value = 202206.6
$ws.Range("L130").Value = 606619.8
$ws.Range("N130").Value = -616659.8
# Row 131
$ws.Range("H131").Value = 1215.7678
$ws.Range("I131").Value = 1132
$ws.Range("J131").Value = 1243.6904
$ws.Range("K131").Value = 3396
$ws.Range("L131").Value = 3731.0712
$ws.Range("M131").Value = 1644
$ws.Range("N131").Value = -13811.0712

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1530.9375
$ws.Range("I113").Value = 1067.7273
$ws.Range("K113").Value = 1067.7273
$ws.Range("M113").Value = 1102.2727

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1762.5
$ws.Range("I16").Value = 2036.3636
$ws.Range("J16").Value = 1160
$ws.Range("K16").Value = 2036.3636
$ws.Range("L16").Value = 1160
$ws.Range("M16").Value = -1866.3636
$ws.Range("N16").Value = -1500
# Row 132
$ws.Range("H132").Value = 6662.5
$ws.Range("I132").Value = 2223.5925
$ws.Range("J132").Value = 14652.533
$ws.Range("K132").Value = 6670.7775
$ws.Range("L132").Value = 43957.599
$ws.Range("M132").Value = -4140.7775
$ws.Range("N132").Value = -49017.599

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 15798777
$ws.Range("I136").Value = 20855592
$ws.Range("J136").Value = 628331.56000000001
$ws.Range("K136").Value = 62566776
$ws.Range("L136").Value = 1884994.68
$ws.Range("M136").Value = -62564226
$ws.Range("N136").Value = -1890094.68
